# Generate Report for Handoff
# Moves the status from "In Translation" to "Ready for handoff" and bumps
# the associated generate/handoff timestamps, across the Overview, zh-cn
# and de-de sheets of the localization status report.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet: per-language status + latest HO xliff generate date ---
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-28 07:04:33"

# --- zh-cn sheet: status + latest handoff datetime ---
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-28 07:04:29"

# --- de-de sheet: status + latest handoff datetime ---
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-28 07:04:33"

# --- Widen the status columns so the longer "Ready for handoff" text fits ---
$wsOverview.Columns("E").ColumnWidth = 16.333333333333332
$wsOverview.Columns("F").ColumnWidth = 16.333333333333332
$wsZhCn.Columns("C").ColumnWidth = 16.333333333333332
$wsDeDe.Columns("C").ColumnWidth = 16.333333333333332
